$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B2 value from 1204 to 120420
$ws.Range("B2").Value = 120420

# Update selection to B2
$ws.Range("B2").Select()
